$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial for every data row (rows 2-370).
# All of these cells currently store 45205 (2023-10-06) and must become
# 45206 (2023-10-07).
$range = $ws.Range("C2:C370")
$range.Value = 45206

